$d = $word.ActiveDocument

# 1. Update "Last update" timestamp
$d.Content.Find.Execute(
    "Last update: 2017-07-24T19:06:50Z", $true, $false, $false, $false, $false,
    $true, 1, $false, "Last update: 2017-07-24T20:54:05Z", 2)

# 2. "location:" -> "geo-location:"
$d.Content.Find.Execute(
    "location:", $true, $false, $false, $false, $false,
    $true, 1, $false, "geo-location:", 2)

# 3. "www-linkedin:" -> "linkedin:"
$d.Content.Find.Execute(
    "www-linkedin:", $true, $false, $false, $false, $false,
    $true, 1, $false, "linkedin:", 2)

# 4. Iskar Matkash paragraph rewording
$d.Content.Find.Execute(
    "Pallets can be automatically move from one stand on a workstation to the next, or to temporarily place a pallet on a storage stand or a stack stand.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The product or intermediate product is placed on pallets. The pallets are moved move from one stand on a workstation to a stand on another workstation, or temporarily to a storage stand or stack.",
    2)

# 5. "We created an implementation architecture" -> "We created the implementation architecture"
$d.Content.Find.Execute(
    "We created an implementation architecture in Pascal.", $true, $false, $false, $false, $false,
    $true, 1, $false, "We created the implementation architecture in Pascal.", 2)

# 6. "4 man-years later" -> "6 man-years later"
$d.Content.Find.Execute(
    "After 18 calendar months, and 4 man-years later, the factory ran perfectly!", $true, $false, $false, $false, $false,
    $true, 1, $false, "After 18 calendar months, and 6 man-years later, the factory ran perfectly!", 2)

# 7. "exhorbitantly priced" -> "exhorbitantly expensive"
$d.Content.Find.Execute(
    "and exhorbitantly priced.", $true, $false, $false, $false, $false,
    $true, 1, $false, "and exhorbitantly expensive.", 2)

# 8. Footnote text rewording
$d.Footnotes.Item(1).Range.Find.Execute(
    "then it is trivial for anyone who may be suspicious of their claim to ask me for verification by sending me an email.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "then it is trivial for anyone who may be suspicious of their claim to " + [char]8220 + "challenge" + [char]8221 + " the site by asking me for a verification email.",
    2)
